# stateless entities outside the US
# The income table gains two new "IMF (20%)" measure columns (Sales / Sales+Emp),
# which are inserted in the column order right before the existing "IMF" columns.
# The old "OECD (20%)" columns are dropped, and the data that used to live in the
# "IMF - Sales" / "IMF - Sales + Emp" columns (F,G) shifts right into the columns
# that used to hold the "OECD (20%)" data (H,I). Columns F,G are then populated
# with the brand new "IMF (20%)" figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 4, 5, 6)

# New "IMF (20%) - Sales" / "IMF (20%) - Sales + Emp" values, keyed by row.
$newF = @{
    2 = 0.7415338153038805
    4 = -0.2969016773552038
    5 = -10.74605673108959
    6 = 0.2028407217499605
}
$newG = @{
    2 = 0.6034852546579005
    4 = 1.497156856423381
    5 = -13.08112174461001
    6 = 0.5157458717548939
}

# Capture the current "IMF - Sales" / "IMF - Sales + Emp" values (columns F, G)
# before overwriting them -- these move into columns H, I.
$oldF = @{}
$oldG = @{}
foreach ($r in $rows) {
    $oldF[$r] = $ws.Cells.Item($r, 6).Value()
    $oldG[$r] = $ws.Cells.Item($r, 7).Value()
}

# Move the old IMF - Sales / IMF - Sales + Emp data into H, I (replacing the old
# OECD (20%) - Sales / OECD (20%) - Sales + Emp data that lived there).
foreach ($r in $rows) {
    $ws.Cells.Item($r, 8).Value = $oldF[$r]
    $ws.Cells.Item($r, 9).Value = $oldG[$r]
}

# Write the new IMF (20%) figures into F, G.
foreach ($r in $rows) {
    $ws.Cells.Item($r, 6).Value = $newF[$r]
    $ws.Cells.Item($r, 7).Value = $newG[$r]
}

# Update the header row so the columns carry their new labels.
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"
$ws.Range("H1").Value = "IMF - Sales"
$ws.Range("I1").Value = "IMF - Sales + Emp"
